$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-18 represent attendance dates; column H is the "Absent" count.
# For most dates the absence flag moves from 0 to 1 (student marked absent).
$absentRows = @(3, 4, 5, 6, 7, 8, 9, 12, 14, 16, 17, 18)
foreach ($r in $absentRows) {
    $ws.Cells.Item($r, 8).Value = 1
}

# For a few dates (rows 10, 11, 13, 15) the "Total Attendance Count" (D)
# and "Real" (E) columns move from 0 to 1 instead, while H stays 0.
$realRows = @(10, 11, 13, 15)
foreach ($r in $realRows) {
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 1
}
